$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "thumbnail" -> "image"
#    This word lives in its own run (distinct rPr/color) inside a paragraph
#    that also has "${" and "}" runs. Rebuild the whole paragraph via
#    Range.InsertXML so every run (incl. the untouched siblings) keeps its
#    exact original formatting/attributes.
# ---------------------------------------------------------------------------
$rngThumb = $d.Content
$foundThumb = $rngThumb.Find.Execute('thumbnail', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($foundThumb) {
    $xmlThumb = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="2E4D27ED" w14:textId="403FD21F" w:rsidR="009C5646" w:rsidRPr="00895CD6" w:rsidRDefault="009C5646" w:rsidP="00CC4E98">
<w:pPr><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00895CD6"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>${</w:t></w:r>
<w:r w:rsidRPr="00895CD6"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:eastAsia="Times New Roman" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:color w:val="CE9178"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="zh-CN"/></w:rPr><w:t>image</w:t></w:r>
<w:r w:rsidRPr="00895CD6"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $rngThumb.InsertXML($xmlThumb)
}

# ---------------------------------------------------------------------------
# 2) "Registration             :${registrationno}" is currently one run;
#    split it into three runs: "Registration             :${",
#    "registrationno", "}" (all sharing identical run formatting).
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute('Registration             :${registrationno}', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) {
    $xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="725CF883" w14:textId="01F1F9FC" w:rsidR="009C5646" w:rsidRPr="00895CD6" w:rsidRDefault="009C5646" w:rsidP="006316D4">
<w:pPr><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00895CD6"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Registration             :${</w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>registrationno</w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $rng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 3) "Registration Date    :${registrationdate}" is currently one run;
#    split it into three runs: "Registration Date    :${",
#    "registrationdate", "}" (all sharing identical run formatting).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute('Registration Date    :${registrationdate}', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($found2) {
    $xml2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="4745F927" w14:textId="396795B2" w:rsidR="009C5646" w:rsidRPr="00895CD6" w:rsidRDefault="009C5646" w:rsidP="006316D4">
<w:pPr><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>
<w:r w:rsidRPr="00895CD6"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Registration Date    :${</w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>registrationdate</w:t></w:r>
<w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
    $rng2.InsertXML($xml2)
}
